# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "Strike#" values previously stored in column G (header "K") are being
# regenerated from the real per-game strikeout counts (K). This script writes
# the corrected K values into column G for every data row on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new K (strikeouts) value for column G.
# Rows not listed here (16, 57, 69, 72) are unchanged by this regeneration.
$kValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 2
    6  = 1
    7  = 3
    8  = 0
    9  = 1
    10 = 2
    11 = 2
    12 = 1
    13 = 2
    14 = 1
    15 = 4
    17 = 2
    18 = 1
    19 = 2
    20 = 1
    21 = 3
    22 = 1
    23 = 2
    24 = 1
    25 = 2
    26 = 0
    27 = 1
    28 = 1
    29 = 3
    30 = 0
    31 = 1
    32 = 2
    33 = 2
    34 = 2
    35 = 1
    36 = 2
    37 = 3
    38 = 2
    39 = 2
    40 = 3
    41 = 0
    42 = 1
    43 = 0
    44 = 2
    45 = 3
    46 = 1
    47 = 2
    48 = 1
    49 = 2
    50 = 4
    51 = 2
    52 = 1
    53 = 3
    54 = 0
    55 = 2
    56 = 2
    58 = 2
    59 = 1
    60 = 1
    61 = 2
    62 = 1
    63 = 1
    64 = 3
    65 = 1
    66 = 2
    67 = 3
    68 = 2
    70 = 2
    71 = 1
    73 = 0
    74 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}

Write-Host "Updated $($kValues.Keys.Count) K (column G) values"
